$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "Used A" column (D) and "Used B" column (originally H, now col 7 after first delete)
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(7).Delete()

# Update "A_sub extra" (col C) values for AMX, CHL, COL rows
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 3

# Update "Bi extra" (col E) values for FOS, POL, TET rows
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 8

# Re-point the "total vol" formulas (col F, shifted from old G) at the new D/E columns
$ws.Range("F2").Formula = "=(B2+C2)*125*8+(D2+E2)*100*6"
$ws.Range("F3").Formula = "=(B3+C3)*125*8+(D3+E3)*100*6"
$ws.Range("F4").Formula = "=(B4+C4)*125*8+(D4+E4)*100*6"
$ws.Range("F5").Formula = "=(B5+C5)*125*8+(D5+E5)*100*6"
$ws.Range("F6").Formula = "=(B6+C6)*125*8+(D6+E6)*100*6"
$ws.Range("F7").Formula = "=(B7+C7)*125*8+(D7+E7)*100*6"

$ws.Range("E8").Select()
